# Table A4 - Descriptives ADES: adapting figure codes and colors
#
# 1) Remove the redundant "Eustress/Distress" sub-header row (the second
#    of the three stacked header rows) - its content duplicated the top
#    header row and is no longer needed.
# 2) Re-code the CL (evaluation) indicator from "0.00"/"1.00" to "0"/"1".
# 3) Swap the "high"/"low" level (and its associated M/SD statistics)
#    between the two rows that share the same CL code, for both the
#    0-code pair and the 1-code pair.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- 1) delete the extra header row (row 2 of 7) ---------------------
$t.Rows.Item(2).Delete()

# After the delete, the table has 6 rows; the four data rows are now
# rows 3-6 (values below reflect the ORIGINAL content of each row,
# before this script touches anything, so they are unambiguous).

# Row 3 was (0.00, high, 2.87, 0.88, 2.34, 0.80) -> becomes the old row 4's
# level/stats, keeping its own (reformatted) CL code.
$t.Cell(3, 1).Range.Text = "0"
$t.Cell(3, 2).Range.Text = "low"
$t.Cell(3, 3).Range.Text = "3.33"
$t.Cell(3, 4).Range.Text = "0.85"
$t.Cell(3, 5).Range.Text = "1.84"
$t.Cell(3, 6).Range.Text = "0.69"

# Row 4 was (0.00, low, 3.33, 0.85, 1.84, 0.69) -> becomes the old row 3's
# level/stats, keeping its own (reformatted) CL code.
$t.Cell(4, 1).Range.Text = "0"
$t.Cell(4, 2).Range.Text = "high"
$t.Cell(4, 3).Range.Text = "2.87"
$t.Cell(4, 4).Range.Text = "0.88"
$t.Cell(4, 5).Range.Text = "2.34"
$t.Cell(4, 6).Range.Text = "0.80"

# Row 5 was (1.00, high, 2.94, 0.97, 2.59, 1.00) -> becomes the old row 6's
# level/stats, keeping its own (reformatted) CL code.
$t.Cell(5, 1).Range.Text = "1"
$t.Cell(5, 2).Range.Text = "low"
$t.Cell(5, 3).Range.Text = "3.17"
$t.Cell(5, 4).Range.Text = "0.91"
$t.Cell(5, 5).Range.Text = "2.35"
$t.Cell(5, 6).Range.Text = "0.95"

# Row 6 was (1.00, low, 3.17, 0.91, 2.35, 0.95) -> becomes the old row 5's
# level/stats, keeping its own (reformatted) CL code.
$t.Cell(6, 1).Range.Text = "1"
$t.Cell(6, 2).Range.Text = "high"
$t.Cell(6, 3).Range.Text = "2.94"
$t.Cell(6, 4).Range.Text = "0.97"
$t.Cell(6, 5).Range.Text = "2.59"
$t.Cell(6, 6).Range.Text = "1.00"

Write-Host "done"
